$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# --- Change 1: Title paragraph "Texas Holdem Poker" -> wrap "Holdem" with proofErr spell tags
$p1 = Find-ParagraphByText $d "Texas Holdem Poker"
$xml1 = '<w:p ' + $wNs + '>' + `
    '<w:r><w:t xml:space="preserve">Texas </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Holdem</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> Poker</w:t></w:r>' + `
    '</w:p>'
$p1.Range.InsertXML($xml1)

# --- Change 2: Body paragraph mentioning "version of Texas Holdem poker." -> same proofErr wrap
$p2 = Find-ParagraphByText $d "In this project, we implemented a"
$xml2 = '<w:p ' + $wNs + '>' + `
    '<w:r><w:t>In this project, we implemented a</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> version of Texas </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Holdem</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> poker. </w:t></w:r>' + `
    '<w:r><w:t>This involved creating a GUI allowing use</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">rs to play poker against a bot. </w:t></w:r>' + `
    '<w:r><w:t>We wrote the code for this project in Python and relied</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> on </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">object oriented programming, specifically </w:t></w:r>' + `
    '<w:r><w:t>inheritance.</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> Our bot AI uses</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> a variatio</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">n of the Monte Carlo </w:t></w:r>' + `
    '<w:r><w:t>Algorithm;</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">the bot </w:t></w:r>' + `
    '<w:r><w:t>uses</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> this algorithm to determine its probability of winning before</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>choosing</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>which</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> move</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> to make.</w:t></w:r>' + `
    '</w:p>'
$p2.Range.InsertXML($xml2)

# --- Change 3: "About Me" bio paragraph rewritten, bookmark moved before the single new run
$p3 = Find-ParagraphByText $d "I am passionate about everything I do"
if (-not $p3) {
    $p3 = Find-ParagraphByText $d "third year of studying Computer Engineering"
}
$newBio = 'I am currently in my third year of studying Computer Engineering at the University of Alberta. I chose this degree for the opportunity to innovate and design solutions to the problems our society faces. I greatly value education and I am in constant pursuit of developing my knowledge and understanding of the industry as well as lifelong learning. I am passionate about everything I do. My greatest goal is to apply this passion to my work and someday improve the lives of others and benefit society.'
$xml3 = '<w:p ' + $wNs + '>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>' + $newBio + '</w:t></w:r>' + `
    '</w:p>'
$p3.Range.InsertXML($xml3)
